# Updates the cryptos list: refresh Price (col D) and Volume(1h) (col E)
# values per-row, and fix the PancakeSwap/RenderToken row order (rows 29/30
# had swapped identities along with their Price/Volume figures).
#
# NOTE: this interpreter's function/param-block named-argument binding is
# unreliable, so every cell is written with a direct, positional statement.
#
# NOTE: the Price column holds text (e.g. "64.099.68", "0.0780") - values
# that look like a single plain decimal number (e.g. "585.87") would
# otherwise get auto-converted to a numeric cell by Excel's normal input
# parsing (losing formatting like trailing zeros, e.g. "0.0780"->0.078), so
# those are written with a leading apostrophe to force a text entry, exactly
# as a user typing into Excel would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "64.099.68"
$ws.Range("E2").Value = "  -0.24%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.482.08"
$ws.Range("E3").Value = "  -0.48%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'585.87"
$ws.Range("E5").Value = "  -0.10%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'132.09"
$ws.Range("E6").Value = "  -1.69%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -0.60%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "'7.68"
$ws.Range("E9").Value = "  +5.29%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = "  -1.47%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "'0.387"
$ws.Range("E11").Value = "  +0.29%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "4.075.32"
$ws.Range("E12").Value = "  -0.53%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.00%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  -1.95%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "3.484.47"
$ws.Range("E15").Value = "  -0.49%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "64.148.46"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17 - Avalanche
$ws.Range("D17").Value = "'24.38"
$ws.Range("E17").Value = "  -7.12%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "'9.99"
$ws.Range("E18").Value = "  +0.59%  "

# Row 19 - Polkadot
$ws.Range("D19").Value = "'5.72"
$ws.Range("E19").Value = "  -0.21%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'13.50"
$ws.Range("E20").Value = "  -2.04%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'384.96"
$ws.Range("E21").Value = "  -1.76%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "'0.575"
$ws.Range("E22").Value = "  +0.59%  "

# Row 23 - WrappedeETH
$ws.Range("D23").Value = "3.626.75"
$ws.Range("E23").Value = "  -0.35%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'74.75"
$ws.Range("E24").Value = "  +0.69%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - LEO
$ws.Range("E26").Value = "  +0.77%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  -1.58%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("E28").Value = "  +0.24%  "

# Rows 29/30 swapped identity (PancakeSwap <-> RenderToken) plus new figures.
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  -4.94%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.23"
$ws.Range("E30").Value = "  -0.05%  "

# Row 31 - Fetch.AI
$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = "  -5.34%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'7.93"
$ws.Range("E32").Value = "  -4.33%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +2.66%  "

# Row 34 - RenzoRestakedETH
$ws.Range("D34").Value = "3.516.96"
$ws.Range("E34").Value = "  -0.11%  "

# Row 35 - USDe (unchanged)

# Row 36 - EthereumClassic
$ws.Range("D36").Value = "'22.96"
$ws.Range("E36").Value = "  -2.25%  "

# Row 37 - NEARProtocol
$ws.Range("D37").Value = "'5.18"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38 - Aptos
$ws.Range("D38").Value = "'6.81"
$ws.Range("E38").Value = "  -1.23%  "

# Row 39 - ImmutableX
$ws.Range("D39").Value = "'1.52"
$ws.Range("E39").Value = "  -3.16%  "

# Row 40 - Monero
$ws.Range("D40").Value = "'163.27"
$ws.Range("E40").Value = "  +0.68%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "'0.0780"
$ws.Range("E41").Value = "  -0.37%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "'0.801"
$ws.Range("E42").Value = "  -0.56%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.13%  "

# Row 44 - Filecoin
$ws.Range("D44").Value = "'4.33"
$ws.Range("E44").Value = "  -1.59%  "

# Row 45 - EnergySwap
$ws.Range("D45").Value = "'24.04"
$ws.Range("E45").Value = "  -5.42%  "

# Row 46 - Stacks
$ws.Range("D46").Value = "'1.62"
$ws.Range("E46").Value = "  -1.66%  "

# Row 47 - ONDO
$ws.Range("D47").Value = "'1.14"
$ws.Range("E47").Value = "  -2.84%  "

# Row 48 - SuiNetwork
$ws.Range("D48").Value = "'0.927"
$ws.Range("E48").Value = "  +3.51%  "

# Row 49 - Cosmos
$ws.Range("D49").Value = "'6.73"
$ws.Range("E49").Value = "  -0.65%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.366.65"
$ws.Range("E50").Value = "  -3.92%  "

# Row 51 - VeChain
$ws.Range("D51").Value = "'0.0255"
$ws.Range("E51").Value = "  -2.54%  "
